# Applies the "Fixed bug that was preventing imports from paste operation" edit
# to DataExample.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text fixes (row 2) ---
# B2: time-format string loses one digit of sub-second precision (OS3 -> OS2)
$ws.Cells.Item(2, 2).Value = "%m/%d/%Y %H:%M:%OS2"
# E2: numeric format string %.3f -> %.2f (this also removes the now-unused
# "%.3f" shared string from the workbook)
$ws.Cells.Item(2, 5).Value = "%.2f"

# --- 2. Column C width grows slightly to fit the new example value ---
$ws.Columns.Item(3).ColumnWidth = 12.8

# --- 3. Fix the mis-styled "many-decimals" demo cell ---
# E21 had accidentally been given the oddball 0.000000-style numeric format;
# put it back to the normal column-E format (0.000, right aligned) ...
$ws.Cells.Item(21, 5).NumberFormat = "0.000"

# ... and move that demonstration format (now with many more decimal places,
# and without the stray right-alignment) onto C7 instead, together with a
# pi-like example value.
$ws.Cells.Item(7, 3).NumberFormat = "0.00000000000"
$ws.Cells.Item(7, 3).Value = 3.14159265358979

# --- 4. A handful of example values in column H change scale ---
$ws.Cells.Item(9, 8).Value = 0.00043
$ws.Cells.Item(13, 8).Value = 0.0029
$ws.Cells.Item(27, 8).Value = -0.000000754

# --- 5. Column M (date) values: rows 9-19 and 21-54 become a consecutive
# run of dates instead of all being identical ---
$ws.Cells.Item(9, 13).Value = 39681
$ws.Cells.Item(10, 13).Value = 39682
$ws.Cells.Item(11, 13).Value = 39683
$ws.Cells.Item(12, 13).Value = 39684
$ws.Cells.Item(13, 13).Value = 39685
$ws.Cells.Item(14, 13).Value = 39686
$ws.Cells.Item(15, 13).Value = 39687
$ws.Cells.Item(16, 13).Value = 39688
$ws.Cells.Item(17, 13).Value = 39689
$ws.Cells.Item(18, 13).Value = 39690
$ws.Cells.Item(19, 13).Value = 39691
$ws.Cells.Item(21, 13).Value = 39692
$ws.Cells.Item(22, 13).Value = 39693
$ws.Cells.Item(23, 13).Value = 39694
$ws.Cells.Item(24, 13).Value = 39695
$ws.Cells.Item(25, 13).Value = 39696
$ws.Cells.Item(26, 13).Value = 39697
$ws.Cells.Item(27, 13).Value = 39698
$ws.Cells.Item(28, 13).Value = 39699
$ws.Cells.Item(29, 13).Value = 39700
$ws.Cells.Item(30, 13).Value = 39701
$ws.Cells.Item(31, 13).Value = 39702
$ws.Cells.Item(32, 13).Value = 39703
$ws.Cells.Item(33, 13).Value = 39704
$ws.Cells.Item(34, 13).Value = 39705
$ws.Cells.Item(35, 13).Value = 39706
$ws.Cells.Item(36, 13).Value = 39707
$ws.Cells.Item(37, 13).Value = 39708
$ws.Cells.Item(38, 13).Value = 39709
$ws.Cells.Item(39, 13).Value = 39710
$ws.Cells.Item(40, 13).Value = 39711
$ws.Cells.Item(41, 13).Value = 39712
$ws.Cells.Item(42, 13).Value = 39713
$ws.Cells.Item(43, 13).Value = 39714
$ws.Cells.Item(44, 13).Value = 39715
$ws.Cells.Item(45, 13).Value = 39716
$ws.Cells.Item(46, 13).Value = 39717
$ws.Cells.Item(47, 13).Value = 39718
$ws.Cells.Item(48, 13).Value = 39719
$ws.Cells.Item(49, 13).Value = 39720
$ws.Cells.Item(50, 13).Value = 39721
$ws.Cells.Item(51, 13).Value = 39722
$ws.Cells.Item(52, 13).Value = 39723
$ws.Cells.Item(53, 13).Value = 39724
$ws.Cells.Item(54, 13).Value = 39725
